# Update the "plan ip" workbook: change the IP address column from the
# 33.x.x.x range to the 172.x.x.x range (same host/subnet suffixes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds the "ip" values for rows 3-9. Replace the "33." network
# prefix with "172." while keeping the rest of each address unchanged.
$ws.Range("D3").Value = "172.0.0.0"
$ws.Range("D4").Value = "172.18.0.0"
$ws.Range("D5").Value = "172.13.0.0"
$ws.Range("D6").Value = "172.18.29.0"
$ws.Range("D7").Value = "172.18.30.0"
$ws.Range("D8").Value = "172.13.28.0"
$ws.Range("D9").Value = "172.13.20.0"

# Update the active selection to match the edited workbook (cursor left on D9).
$ws.Range("D9").Select()
